$wb = $excel.ActiveWorkbook

# Delete the "1486" (LİHKAB) record row from the "Kayitlar" sheet.
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(619).Delete()

# Delete the same record row from the "Merkez İlçe" sheet (its filtered view).
$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(77).Delete()
